$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.99788412191014
$ws.Range("C2").Value = 9.326518530489631
$ws.Range("E2").Value = 14.41306444996972
$ws.Range("F2").Value = 39.25783281856497
$ws.Range("G2").Value = 39.3468876080474
$ws.Range("H2").Value = 16.76722304247926
$ws.Range("I2").Value = 26.04045985153383
$ws.Range("J2").Value = 8.150861607629425
$ws.Range("L2").Value = 12.48407712883907
$ws.Range("M2").Value = 17.25861238814616
$ws.Range("N2").Value = 19.05130294312527
$ws.Range("B3").Value = 16.59957838986409
$ws.Range("C3").Value = 8.983752758379715
$ws.Range("E3").Value = 14.42940836260917
$ws.Range("F3").Value = 39.26928952732694
$ws.Range("G3").Value = 39.26067236931777
$ws.Range("H3").Value = 16.806263338886
$ws.Range("I3").Value = 26.12642817044946
$ws.Range("J3").Value = 8.157720764025218
$ws.Range("L3").Value = 12.48674129313509
$ws.Range("M3").Value = 17.18251165123601
$ws.Range("N3").Value = 19.11705770635037
$ws.Range("B4").Value = 16.35350579027716
$ws.Range("C4").Value = 8.764543140776011
$ws.Range("E4").Value = 14.4401641084054
$ws.Range("F4").Value = 39.28576578983733
$ws.Range("G4").Value = 39.22121054910279
$ws.Range("H4").Value = 16.83368068097822
$ws.Range("I4").Value = 26.18484011123047
$ws.Range("J4").Value = 8.162182109599648
$ws.Range("L4").Value = 12.48982704942565
$ws.Range("M4").Value = 17.13834541705754
$ws.Range("N4").Value = 19.15930016257053
$ws.Range("B5").Value = 16.25300017608063
$ws.Range("C5").Value = 8.673086493014198
$ws.Range("E5").Value = 14.44472876194452
$ws.Range("F5").Value = 39.29485279665398
$ws.Range("G5").Value = 39.20852421116899
$ws.Range("H5").Value = 16.8457182044415
$ws.Range("I5").Value = 26.21005464577235
$ws.Range("J5").Value = 8.164063200236486
$ws.Range("L5").Value = 12.49145015741187
$ws.Range("M5").Value = 17.12100410441971
$ws.Range("N5").Value = 19.17698572607394
$ws.Range("B6").Value = 16.23630198021537
$ws.Range("C6").Value = 8.657774044300737
$ws.Range("E6").Value = 14.44549769975837
$ws.Range("F6").Value = 39.29650493164932
$ws.Range("G6").Value = 39.20662278606976
$ws.Range("H6").Value = 16.84776919570395
$ws.Range("I6").Value = 26.2143266206559
$ws.Range("J6").Value = 8.164379369471623
$ws.Range("L6").Value = 12.49174178554803
$ws.Range("M6").Value = 17.11816462793683
$ws.Range("N6").Value = 19.17995091725231
$ws.Range("B7").Value = 16.35215106029798
$ws.Range("C7").Value = 8.763318232841479
$ws.Range("E7").Value = 14.44022493305425
$ws.Range("F7").Value = 39.28587873590683
$ws.Range("G7").Value = 39.22102570668173
$ws.Range("H7").Value = 16.83383952459916
$ws.Range("I7").Value = 26.18517445448774
$ws.Range("J7").Value = 8.162207223019873
$ws.Range("L7").Value = 12.48984745746405
$ws.Range("M7").Value = 17.13810887006143
$ws.Range("N7").Value = 19.1595367655264
$ws.Range("B8").Value = 16.86095079976051
$ws.Range("C8").Value = 9.210190200593503
$ws.Range("E8").Value = 14.41855053957188
$ws.Range("F8").Value = 39.25982296360402
$ws.Range("G8").Value = 39.31436902470318
$ws.Range("H8").Value = 16.77996772084293
$ws.Range("I8").Value = 26.06893159409484
$ws.Range("J8").Value = 8.153174945754772
$ws.Range("L8").Value = 12.48469543854892
$ws.Range("M8").Value = 17.23184960731734
$ws.Range("N8").Value = 19.07358823130633
$ws.Range("B9").Value = 17.84038741821795
$ws.Range("C9").Value = 10.01429796502007
$ws.Range("E9").Value = 14.38174544469971
$ws.Range("F9").Value = 39.28365625284484
$ws.Range("G9").Value = 39.60385368565238
$ws.Range("H9").Value = 16.70175191449395
$ws.Range("I9").Value = 25.88579332059169
$ws.Range("J9").Value = 8.137433759458286
$ws.Range("L9").Value = 12.48604953813173
$ws.Range("M9").Value = 17.43538025697428
$ws.Range("N9").Value = 18.91980170636602
$ws.Range("B10").Value = 18.54083903645992
$ws.Range("C10").Value = 10.55789859906869
$ws.Range("E10").Value = 14.35815319628381
$ws.Range("F10").Value = 39.34678633868141
$ws.Range("G10").Value = 39.88050224013845
$ws.Range("H10").Value = 16.66111043582221
$ws.Range("I10").Value = 25.77878288758604
$ws.Range("J10").Value = 8.127055235406521
$ws.Range("L10").Value = 12.49396226555705
$ws.Range("M10").Value = 17.59607092585657
$ws.Range("N10").Value = 18.81571314112175
$ws.Range("B11").Value = 18.85380821658132
$ws.Range("C11").Value = 10.79442722627318
$ws.Range("E11").Value = 14.34816393995575
$ws.Range("F11").Value = 39.38537093770116
$ws.Range("G11").Value = 40.01996003279721
$ws.Range("H11").Value = 16.64629391547804
$ws.Range("I11").Value = 25.7361267836072
$ws.Range("J11").Value = 8.122588210339174
$ws.Range("L11").Value = 12.49904823190547
$ws.Range("M11").Value = 17.6714005637889
$ws.Range("N11").Value = 18.77027173414248
$ws.Range("B12").Value = 18.97139043096632
$ws.Range("C12").Value = 10.88241155790651
$ws.Range("E12").Value = 14.34448769549276
$ws.Range("F12").Value = 39.40139476425674
$ws.Range("G12").Value = 40.07469546068037
$ws.Range("H12").Value = 16.64121253232673
$ws.Range("I12").Value = 25.72084375129933
$ws.Range("J12").Value = 8.120932970234847
$ws.Range("L12").Value = 12.50118637876556
$ws.Range("M12").Value = 17.70022935739118
$ws.Range("N12").Value = 18.75333721007906
$ws.Range("B13").Value = 18.94611023800531
$ws.Range("C13").Value = 10.86353353446568
$ws.Range("E13").Value = 14.34527471106179
$ws.Range("F13").Value = 39.39788102979367
$ws.Range("G13").Value = 40.06282209499992
$ws.Range("H13").Value = 16.64228333811387
$ws.Range("I13").Value = 25.72409647828863
$ws.Range("J13").Value = 8.121287843888844
$ws.Range("L13").Value = 12.50071647842006
$ws.Range("M13").Value = 17.69400736682607
$ws.Range("N13").Value = 18.75697223592515
$ws.Range("B14").Value = 18.86350110272921
$ws.Range("C14").Value = 10.80169768685433
$ws.Range("E14").Value = 14.34785936138885
$ws.Range("F14").Value = 39.38666095548231
$ws.Range("G14").Value = 40.02442474622737
$ws.Range("H14").Value = 16.64586525169886
$ws.Range("I14").Value = 25.73485198416254
$ws.Range("J14").Value = 8.122451306043963
$ws.Range("L14").Value = 12.49921989486643
$ws.Range("M14").Value = 17.67376635435051
$ws.Range("N14").Value = 18.76887305406522
$ws.Range("B15").Value = 18.81277586790553
$ws.Range("C15").Value = 10.76361416686631
$ws.Range("E15").Value = 14.34945638902128
$ws.Range("F15").Value = 39.37997210250169
$ws.Range("G15").Value = 40.00115506361657
$ws.Range("H15").Value = 16.64812824641087
$ws.Range("I15").Value = 25.74155343841521
$ws.Range("J15").Value = 8.123168684008698
$ws.Range("L15").Value = 12.49833078176216
$ws.Range("M15").Value = 17.66140708938713
$ws.Range("N15").Value = 18.77619818088208
$ws.Range("B16").Value = 18.52026146164314
$ws.Range("C16").Value = 10.54222113562747
$ws.Range("E16").Value = 14.35882093662071
$ws.Range("F16").Value = 39.34446279309122
$ws.Range("G16").Value = 39.87165961938467
$ws.Range("H16").Value = 16.66215272993244
$ws.Range("I16").Value = 25.78169209494235
$ws.Range("J16").Value = 8.127352261746461
$ws.Range("L16").Value = 12.49365966645885
$ws.Range("M16").Value = 17.59119134341622
$ws.Range("N16").Value = 18.81872114124092
$ws.Range("B17").Value = 18.33927506008573
$ws.Range("C17").Value = 10.40362025473369
$ws.Range("E17").Value = 14.36475581335559
$ws.Range("F17").Value = 39.32520212221512
$ws.Range("G17").Value = 39.79568303816487
$ws.Range("H17").Value = 16.67169767321136
$ws.Range("I17").Value = 25.807861322792
$ws.Range("J17").Value = 8.129983697780082
$ws.Range("L17").Value = 12.49117372819921
$ws.Range("M17").Value = 17.5486751590444
$ws.Range("N17").Value = 18.84529556842236
$ws.Range("B18").Value = 18.23464924825191
$ws.Range("C18").Value = 10.32288993551835
$ws.Range("E18").Value = 14.3682393517227
$ws.Range("F18").Value = 39.31505302317282
$ws.Range("G18").Value = 39.75326638345796
$ws.Range("H18").Value = 16.67753318005476
$ws.Range("I18").Value = 25.82347991358468
$ws.Range("J18").Value = 8.131521172042195
$ws.Range("L18").Value = 12.48988389002272
$ws.Range("M18").Value = 17.52443202256431
$ws.Range("N18").Value = 18.86076024098012
$ws.Range("B19").Value = 18.19913790830052
$ws.Range("C19").Value = 10.29538365345936
$ws.Range("E19").Value = 14.36943084396087
$ws.Range("F19").Value = 39.31177645383775
$ws.Range("G19").Value = 39.73912608511635
$ws.Range("H19").Value = 16.67956827961086
$ws.Range("I19").Value = 25.8288653225737
$ws.Range("J19").Value = 8.132045853509437
$ws.Range("L19").Value = 12.48947126500792
$ws.Range("M19").Value = 17.51626049150585
$ws.Range("N19").Value = 18.86602722992911
$ws.Range("B20").Value = 18.35859681315188
$ws.Range("C20").Value = 10.41847944965682
$ws.Range("E20").Value = 14.36411679829363
$ws.Range("F20").Value = 39.32715633167751
$ws.Range("G20").Value = 39.8036382955914
$ws.Range("H20").Value = 16.67064582852458
$ws.Range("I20").Value = 25.80501688045311
$ws.Range("J20").Value = 8.129701100697922
$ws.Range("L20").Value = 12.49142388119494
$ws.Range("M20").Value = 17.55317936325164
$ws.Range("N20").Value = 18.84244808017499
$ws.Range("B21").Value = 18.88779160352884
$ws.Range("C21").Value = 10.81990362918743
$ws.Range("E21").Value = 14.34709730028695
$ws.Range("F21").Value = 39.38991827941273
$ws.Range("G21").Value = 40.03565098086473
$ws.Range("H21").Value = 16.644798781556
$ws.Range("I21").Value = 25.73166919319736
$ws.Range("J21").Value = 8.122108584980284
$ws.Range("L21").Value = 12.49965373140151
$ws.Range("M21").Value = 17.67970354753882
$ws.Range("N21").Value = 18.76537009285449
$ws.Range("B22").Value = 19.22816801208262
$ws.Range("C22").Value = 11.07301175046942
$ws.Range("E22").Value = 14.33659449554841
$ws.Range("F22").Value = 39.43916715443557
$ws.Range("G22").Value = 40.19849246096642
$ws.Range("H22").Value = 16.63099197854527
$ws.Range("I22").Value = 25.68880436565274
$ws.Range("J22").Value = 8.117358062683806
$ws.Range("L22").Value = 12.50626850767807
$ws.Range("M22").Value = 17.76415258005084
$ws.Range("N22").Value = 18.71658668056246
$ws.Range("B23").Value = 19.04704082466001
$ws.Range("C23").Value = 10.93878000926868
$ws.Range("E23").Value = 14.34214339231061
$ws.Range("F23").Value = 39.41213141585317
$ws.Range("G23").Value = 40.11056668911813
$ws.Range("H23").Value = 16.63807818198607
$ws.Range("I23").Value = 25.71121683236579
$ws.Range("J23").Value = 8.119874218045583
$ws.Range("L23").Value = 12.50262550833804
$ws.Range("M23").Value = 17.71892562181265
$ws.Range("N23").Value = 18.74247812916412
$ws.Range("B24").Value = 18.34986322843043
$ws.Range("C24").Value = 10.41176486552369
$ws.Range("E24").Value = 14.36440547428926
$ws.Range("F24").Value = 39.32626995451496
$ws.Range("G24").Value = 39.80003778622785
$ws.Range("H24").Value = 16.67112028358527
$ws.Range("I24").Value = 25.80630106646502
$ws.Range("J24").Value = 8.129828786095951
$ws.Range("L24").Value = 12.49131035282031
$ws.Range("M24").Value = 17.55114238823872
$ws.Range("N24").Value = 18.84373484797228
$ws.Range("B25").Value = 17.57823832000876
$ws.Range("C25").Value = 9.804850366549166
$ws.Range("E25").Value = 14.39109479379097
$ws.Range("F25").Value = 39.26918767101044
$ws.Range("G25").Value = 39.51422374472403
$ws.Range("H25").Value = 16.71996385341093
$ws.Range("I25").Value = 25.93051614553584
$ws.Range("J25").Value = 8.141482715011319
$ws.Range("L25").Value = 12.48446318637896
$ws.Range("M25").Value = 17.3782960113537
$ws.Range("N25").Value = 18.95983530791342
